$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26..137 down to 27..138
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new weekly data point.
# (Columns that are identical to the old row 26 content are simply re-set
# for completeness / robustness.)
$ws.Cells.Item(26, 1).Value = 11
$ws.Cells.Item(26, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(26, 3).Value = "Bíobío"
$ws.Cells.Item(26, 4).Value = 44971
$ws.Cells.Item(26, 5).Value = 8
$ws.Cells.Item(26, 6).Value = 100112001
$ws.Cells.Item(26, 7).Value = "Berenjena"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 100
$ws.Cells.Item(26, 11).Value = 8500
$ws.Cells.Item(26, 12).Value = 9000
$ws.Cells.Item(26, 13).Value = 8750
$ws.Cells.Item(26, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(26, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(26, 16).Value = 146
$ws.Cells.Item(26, 17).Value = 60
$ws.Cells.Item(26, 18).Value = "Hortaliza"

# Match the date cell's number format style used by the rest of column D
$ws.Cells.Item(26, 4).NumberFormat = $ws.Cells.Item(27, 4).NumberFormat
